# Add season record columns (Wins, Losses, Ties) to the player stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy formatting from an existing header cell (A1)
# so the new headers share the same bold/centered/bordered style (s="1").
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# --- Data rows (rows 2-48): season record is the same for every player row
# in this sheet (team-level Wins/Losses/Ties repeated per player).
$wins = 67
$losses = 95
$ties = 0

for ($row = 2; $row -le 48; $row++) {
    $ws.Cells.Item($row, 29).Value = $wins    # column AC
    $ws.Cells.Item($row, 30).Value = $losses  # column AD
    $ws.Cells.Item($row, 31).Value = $ties    # column AE
}
